$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("Albatros" row) - remaining rows shift up by one automatically
$ws.Rows.Item(2).Delete()
